$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: change "Java" -> "Java code" in D2
$ws.Range("D2").Value = "Java code"

# Row 3: shift columns - C3 becomes old D3 text, D3 becomes "Java code", E3 becomes old E4 text
$ws.Range("C3").Value = "Longest Palindromic Subsequence"
$ws.Range("D3").Value = "Java code"
$ws.Range("E3").Value = "Tabulation+space optimization"

# Row 4
$ws.Range("A4").Value = "GFG"
$ws.Range("B4").Value = "GFG"
$ws.Range("C4").Value = "Longest Common Substring"
$ws.Range("D4").Value = "Java code"
$ws.Range("E4").Value = "Tabulation+space optimization"

# Row 5
$ws.Range("A5").Value = "GFG"
$ws.Range("B5").Value = "GFG"
$ws.Range("C5").Value = "Print all LCS sequences"
$ws.Range("D5").Value = "Java code"
$ws.Range("E5").Value = "NA"

# Row 6 (new content, was previously blank)
$ws.Range("A6").Value = 1312
$ws.Range("B6").Value = "LC"
$ws.Range("C6").Value = "Minimum Insertion Steps to Make a String Palindrome"
$ws.Range("D6").Value = "Java code"
$ws.Range("E6").Value = "space optimization -->LCS "

# Row 7 (new content), row height changes from 30.75 to 13.5
$ws.Range("A7").Value = 583
$ws.Range("B7").Value = "LC"
$ws.Range("C7").Value = "Delete Operation for Two Strings"
$ws.Range("D7").Value = "Java code"
$ws.Range("E7").Value = "space optimization -->LCS "
$ws.Rows.Item(7).RowHeight = 13.5

# Row 8 (new content)
$ws.Range("A8").Value = 1092
$ws.Range("B8").Value = "LC"
$ws.Range("C8").Value = "Shortest Common Supersequence "
$ws.Range("D8").Value = "Java code"
$ws.Range("E8").Value = "space optimization -->LCS "

# Update the active selection to A8
$ws.Range("A8").Select()
